$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-08-20"

# Update the header label in I1 (shared string "2022 (through 08-19)" -> "2022 (through 08-20)")
$ws.Range("I1").Value = "2022 (through 08-20)"

# Update September total (row 9) and grand Total row (row 14) for column I
$ws.Range("I9").Value = 119
$ws.Range("I14").Value = 1090
